# vent.xlsx — n_50 tables: update assumed n_50 values (tested multiple_run_results)
$wb = $excel.ActiveWorkbook

# --- n_50_table_1957 ---------------------------------------------------
$ws1957 = $wb.Worksheets.Item("n_50_table_1957")
$ws1957.Select()

$ws1957.Range("B6:B14").Value = 8
$ws1957.Range("B15:B17").Value = 3
$ws1957.Range("B18").Value = 8

$ws1957.Range("B40").Select()

# --- n_50_table_1978 ---------------------------------------------------
$ws1978 = $wb.Worksheets.Item("n_50_table_1978")
$ws1978.Select()

$ws1978.Range("B4:B5").Value = 11
$ws1978.Range("B22:B23").Value = 7
$ws1978.Range("B24:B32").Value = 3
$ws1978.Range("B33:B35").Value = 1.5
$ws1978.Range("B36").Value = 3

$ws1978.Range("B38").Select()

# --- n_50_table_1994 ---------------------------------------------------
$ws1994 = $wb.Worksheets.Item("n_50_table_1994")
$ws1994.Select()

$ws1994.Range("B6:B14").Value = 3
$ws1994.Range("B18").Value = 3
$ws1994.Range("B33:B35").Value = 1

$ws1994.Range("B38").Select()
